# Journal de bord.docx - apply the two edits described by the diff:
#  1. Merge the three "Diagramme UML ( en cas d'utilisation)" runs (which
#     had proofErr gramStart/gramEnd markers around "( en") into a single
#     run with the full text and no xml:space="preserve".
#  2. Append two new paragraphs at the very end of the document (right
#     before the final sectPr), after the existing trailing empty
#     paragraph: a "Difficulté rencontré :" heading-like line (carrying a
#     lastRenderedPageBreak marker on its run) and a long paragraph of
#     text describing the difficulty encountered.

$d = $word.ActiveDocument

# --- Change 1: merge the UML diagram bullet runs -------------------------
# FindWhat, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace(wdReplaceAll)
[void]$d.Content.Find.Execute(
    "Diagramme UML ( en cas d’utilisation)",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Diagramme UML ( en cas d’utilisation)",
    2)

# --- Change 2: add the "Difficulté rencontré" paragraphs at the end ------
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$lastPara = $d.Paragraphs.Last
[void]$lastPara.Range.InsertParagraphAfter()

$diffParaXml = "<w:p $wNs>" +
    "<w:pPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:lastRenderedPageBreak/><w:t>Difficulté rencontré :</w:t></w:r>" +
    "</w:p>"
[void]$d.Paragraphs.Last.Range.InsertXML($diffParaXml)

$bodyText = "On n’a pas pu envoyez de requête via l’Ethernet Shield suite à des problèmes de réseaux et de câblage, mais surtout de code. On ne pouvait également pas connecter la mannette directement au bras robot via arduino car il aurait fallu ouvrir complètement la manette. "
[void]$d.Paragraphs.Last.Range.InsertParagraphAfter()

$bodyParaXml = "<w:p $wNs>" +
    "<w:pPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>$bodyText</w:t></w:r>" +
    "</w:p>"
[void]$d.Paragraphs.Last.Range.InsertXML($bodyParaXml)
